$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) "244.52"
Set-TextValue $ws.Cells.Item(2, 7) "2"
Set-TextValue $ws.Cells.Item(3, 4) "23.15"
Set-TextValue $ws.Cells.Item(3, 7) "2"
Set-TextValue $ws.Cells.Item(4, 4) "5.413"
Set-TextValue $ws.Cells.Item(4, 7) "2"
Set-TextValue $ws.Cells.Item(5, 4) "0.05988"
Set-TextValue $ws.Cells.Item(5, 7) "2"
Set-TextValue $ws.Cells.Item(6, 4) "3.468"
Set-TextValue $ws.Cells.Item(6, 7) "2"
Set-TextValue $ws.Cells.Item(7, 4) "6.523"
Set-TextValue $ws.Cells.Item(7, 7) "2"
Set-TextValue $ws.Cells.Item(8, 4) "0.8189"
Set-TextValue $ws.Cells.Item(8, 7) "2"
Set-TextValue $ws.Cells.Item(9, 4) "0.9159"
Set-TextValue $ws.Cells.Item(9, 7) "2"
Set-TextValue $ws.Cells.Item(10, 4) "0.1412"
Set-TextValue $ws.Cells.Item(10, 7) "2"
Set-TextValue $ws.Cells.Item(11, 4) "0.07388"
Set-TextValue $ws.Cells.Item(11, 7) "2"
Set-TextValue $ws.Cells.Item(12, 4) "0.03255"
Set-TextValue $ws.Cells.Item(12, 7) "2"
Set-TextValue $ws.Cells.Item(13, 4) "0.03042"
Set-TextValue $ws.Cells.Item(13, 7) "2"
Set-TextValue $ws.Cells.Item(14, 4) "0.09349"
Set-TextValue $ws.Cells.Item(14, 7) "2"
Set-TextValue $ws.Cells.Item(15, 4) "3.858"
Set-TextValue $ws.Cells.Item(15, 7) "2"
Set-TextValue $ws.Cells.Item(16, 7) "2"
Set-TextValue $ws.Cells.Item(17, 4) "0.04687"
Set-TextValue $ws.Cells.Item(17, 7) "2"
Set-TextValue $ws.Cells.Item(18, 4) "0.0005940"
Set-TextValue $ws.Cells.Item(18, 7) "2"
Set-TextValue $ws.Cells.Item(19, 4) "0.006101"
Set-TextValue $ws.Cells.Item(19, 7) "2"
Set-TextValue $ws.Cells.Item(20, 4) "0.005024"
Set-TextValue $ws.Cells.Item(20, 5) "19HotbitTokenHTBBestin24h"
Set-TextValue $ws.Cells.Item(20, 7) "2"
Set-TextValue $ws.Cells.Item(21, 4) "0.0009835"
Set-TextValue $ws.Cells.Item(21, 7) "2"
Set-TextValue $ws.Cells.Item(22, 4) "0.00007899"
Set-TextValue $ws.Cells.Item(22, 7) "2"
Set-TextValue $ws.Cells.Item(23, 4) "0.0002900"
Set-TextValue $ws.Cells.Item(23, 7) "2"
Set-TextValue $ws.Cells.Item(24, 4) "3.634"
Set-TextValue $ws.Cells.Item(24, 7) "2"
Set-TextValue $ws.Cells.Item(25, 4) "2.144"
Set-TextValue $ws.Cells.Item(25, 7) "2"
Set-TextValue $ws.Cells.Item(26, 4) "0.3204"
Set-TextValue $ws.Cells.Item(26, 7) "2"
Set-TextValue $ws.Cells.Item(27, 4) "0.1331"
Set-TextValue $ws.Cells.Item(27, 7) "2"
Set-TextValue $ws.Cells.Item(28, 7) "2"
Set-TextValue $ws.Cells.Item(29, 7) "2"
Set-TextValue $ws.Cells.Item(30, 7) "2"
Set-TextValue $ws.Cells.Item(31, 7) "2"
Set-TextValue $ws.Cells.Item(32, 7) "2"
Set-TextValue $ws.Cells.Item(33, 7) "2"
Set-TextValue $ws.Cells.Item(34, 7) "2"
Set-TextValue $ws.Cells.Item(35, 7) "2"
Set-TextValue $ws.Cells.Item(36, 7) "2"
Set-TextValue $ws.Cells.Item(37, 7) "2"
Set-TextValue $ws.Cells.Item(38, 7) "2"
Set-TextValue $ws.Cells.Item(39, 7) "2"
Set-TextValue $ws.Cells.Item(40, 4) "0.03918"
Set-TextValue $ws.Cells.Item(40, 7) "2"
Set-TextValue $ws.Cells.Item(41, 4) "0.006233"
Set-TextValue $ws.Cells.Item(41, 5) "40KickTokenKICK"
Set-TextValue $ws.Cells.Item(41, 7) "2"
Set-TextValue $ws.Cells.Item(42, 4) "0.1075"
Set-TextValue $ws.Cells.Item(42, 7) "2"
Set-TextValue $ws.Cells.Item(43, 4) "0.002569"
Set-TextValue $ws.Cells.Item(43, 7) "2"
Set-TextValue $ws.Cells.Item(44, 4) "0.006498"
Set-TextValue $ws.Cells.Item(44, 7) "2"
Set-TextValue $ws.Cells.Item(45, 4) "0.00005252"
Set-TextValue $ws.Cells.Item(45, 7) "2"
Set-TextValue $ws.Cells.Item(46, 7) "2"
Set-TextValue $ws.Cells.Item(47, 7) "2"
Set-TextValue $ws.Cells.Item(48, 4) "0.9101"
Set-TextValue $ws.Cells.Item(48, 7) "2"
Set-TextValue $ws.Cells.Item(49, 7) "2"
Set-TextValue $ws.Cells.Item(50, 4) "0.00002100"
Set-TextValue $ws.Cells.Item(50, 7) "2"
Set-TextValue $ws.Cells.Item(51, 7) "2"
